$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First 10 data rows (rows 2-11) already exist; update their contents in place.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "5010408-92.2022.8.21.0014"
$ws.Cells.Item(2, 3).Value = "5002665-02.2020.8.21.0014"
$ws.Cells.Item(2, 4).Value = "CIV.36852.01"
$ws.Cells.Item(2, 5).Value = "originario_principal"

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "5217136-10.2022.8.21.0001"
$ws.Cells.Item(3, 3).Value = "5004829-76.2020.8.21.0001"
$ws.Cells.Item(3, 4).Value = "CIV.35064.01"
$ws.Cells.Item(3, 5).Value = "originario_principal"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "5001809-74.2017.8.21.0036"
$ws.Cells.Item(4, 3).Value = "9002336-55.2017.8.21.0036"
$ws.Cells.Item(4, 4).Value = "CIV.04873.01"
$ws.Cells.Item(4, 5).Value = "originario_principal"

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "5002735-26.2018.8.21.0002"
$ws.Cells.Item(5, 3).Value = "9000271-58.2018.8.21.0002"
$ws.Cells.Item(5, 4).Value = "CIV.13978.01"
$ws.Cells.Item(5, 5).Value = "originario_principal"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "5002731-86.2018.8.21.0002"
$ws.Cells.Item(6, 3).Value = "9000433-53.2018.8.21.0002"
$ws.Cells.Item(6, 4).Value = "CIV.36875.01"
$ws.Cells.Item(6, 5).Value = "originario_principal"

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "5035108-55.2014.8.21.0001"
$ws.Cells.Item(7, 3).Value = "0028693-44.2014.8.21.0001"
$ws.Cells.Item(7, 4).Value = "CIV.27992.01"
$ws.Cells.Item(7, 5).Value = "originario_principal"

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "5035108-55.2014.8.21.0001"
$ws.Cells.Item(8, 3).Value = "0028693-44.2014.8.21.0001"
$ws.Cells.Item(8, 4).Value = "CIV.27992.01"
$ws.Cells.Item(8, 5).Value = "originario_principal"

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "5035108-55.2014.8.21.0001"
$ws.Cells.Item(9, 3).Value = "0028693-44.2014.8.21.0001"
$ws.Cells.Item(9, 4).Value = "CIV.27992.01"
$ws.Cells.Item(9, 5).Value = "originario_principal"

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "5035108-55.2014.8.21.0001"
$ws.Cells.Item(10, 3).Value = "0028693-44.2014.8.21.0001"
$ws.Cells.Item(10, 4).Value = "CIV.27992.01"
$ws.Cells.Item(10, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "5009658-08.2018.8.21.0022"
$ws.Cells.Item(11, 3).Value = "9003683-34.2018.8.21.0022"
$ws.Cells.Item(11, 4).Value = "CIV.10163.01"
$ws.Cells.Item(11, 5).Value = "originario_principal"

# New rows (12-30): copy formatting from row 11 first (to get the A-column border/bold style),
# then populate the values.
$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(12, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "5002111-39.2019.8.21.0067"
$ws.Cells.Item(12, 3).Value = "9000482-59.2019.8.21.0067"
$ws.Cells.Item(12, 4).Value = "CIV.11871.01"
$ws.Cells.Item(12, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(13, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "5002115-76.2019.8.21.0067"
$ws.Cells.Item(13, 3).Value = "9001255-07.2019.8.21.0067"
$ws.Cells.Item(13, 4).Value = "CIV.14996.01"
$ws.Cells.Item(13, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(14, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "5001005-76.2018.8.21.0067"
$ws.Cells.Item(14, 3).Value = "9000856-12.2018.8.21.0067"
$ws.Cells.Item(14, 4).Value = "CIV.10944.01"
$ws.Cells.Item(14, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "5002061-76.2020.8.21.0067"
$ws.Cells.Item(15, 3).Value = "9000182-63.2020.8.21.0067"
$ws.Cells.Item(15, 4).Value = "CIV.34197.01"
$ws.Cells.Item(15, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "5002116-61.2019.8.21.0067"
$ws.Cells.Item(16, 3).Value = "9000433-18.2019.8.21.0067"
$ws.Cells.Item(16, 4).Value = "CIV.04100.01"
$ws.Cells.Item(16, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "5002062-61.2020.8.21.0067"
$ws.Cells.Item(17, 3).Value = "9000400-91.2020.8.21.0067"
$ws.Cells.Item(17, 4).Value = "CIV.36413.01"
$ws.Cells.Item(17, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "5001006-61.2018.8.21.0067"
$ws.Cells.Item(18, 3).Value = "9000817-15.2018.8.21.0067"
$ws.Cells.Item(18, 4).Value = "CIV.35921.01"
$ws.Cells.Item(18, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(19, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "5002063-46.2020.8.21.0067"
$ws.Cells.Item(19, 3).Value = "9000178-26.2020.8.21.0067"
$ws.Cells.Item(19, 4).Value = "CIV.34187.01"
$ws.Cells.Item(19, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(20, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "5002065-16.2020.8.21.0067"
$ws.Cells.Item(20, 3).Value = "9000388-77.2020.8.21.0067"
$ws.Cells.Item(20, 4).Value = "CIV.36384.01"
$ws.Cells.Item(20, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(21, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "5001007-46.2018.8.21.0067"
$ws.Cells.Item(21, 3).Value = "9000790-32.2018.8.21.0067"
$ws.Cells.Item(21, 4).Value = "CIV.35935.01"
$ws.Cells.Item(21, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(22, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "5002068-68.2020.8.21.0067"
$ws.Cells.Item(22, 3).Value = "9000458-94.2020.8.21.0067"
$ws.Cells.Item(22, 4).Value = "CIV.36605.01"
$ws.Cells.Item(22, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(23, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "5000282-08.2011.8.21.0001"
$ws.Cells.Item(23, 3).Value = "0302109-66.2011.8.21.0001"
$ws.Cells.Item(23, 4).Value = "CIV.43266.01"
$ws.Cells.Item(23, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(24, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "5034992-20.2012.8.21.0001"
$ws.Cells.Item(24, 3).Value = "0303805-06.2012.8.21.0001"
$ws.Cells.Item(24, 4).Value = "CIV.12799.01"
$ws.Cells.Item(24, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(25, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "5008342-82.2022.8.21.0033"
$ws.Cells.Item(25, 3).Value = "9001156-08.2020.8.21.0033"
$ws.Cells.Item(25, 4).Value = "CIV.17400.01"
$ws.Cells.Item(25, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(26, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "5009659-90.2018.8.21.0022"
$ws.Cells.Item(26, 3).Value = "9005551-47.2018.8.21.0022"
$ws.Cells.Item(26, 4).Value = "CIV.08273.01"
$ws.Cells.Item(26, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(27, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "5002741-33.2018.8.21.0002"
$ws.Cells.Item(27, 3).Value = "9000456-96.2018.8.21.0002"
$ws.Cells.Item(27, 4).Value = "CIV.04071.01"
$ws.Cells.Item(27, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(28, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "5008759-20.2011.8.21.0001"
$ws.Cells.Item(28, 3).Value = "0148964-87.2011.8.21.0001"
$ws.Cells.Item(28, 4).Value = "CIV.00159.01"
$ws.Cells.Item(28, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(29, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "5019558-25.2011.8.21.0001"
$ws.Cells.Item(29, 3).Value = "0218063-47.2011.8.21.0001"
$ws.Cells.Item(29, 4).Value = "CIV.29683.01"
$ws.Cells.Item(29, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Copy() | Out-Null
$ws.Cells.Item(30, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "5125809-52.2020.8.21.0001"
$ws.Cells.Item(30, 3).Value = "9001529-57.2020.8.21.0027"
$ws.Cells.Item(30, 4).Value = "CIV.34964.01"
$ws.Cells.Item(30, 5).Value = "originario_principal"

$excel.CutCopyMode = 0